$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the wording of the expected result of test case MODEL_EXEC_1 (F7):
# the output folder is "ARE\data", not "ARE\data\csv"
$ws.Range("F7").Value = "1. The model must be deployed and started successfully.`n2. The Oscilloscope must show a correct sinus signal`n3. After clicking 'Stop Writing' a new file must exist in the folder: ARE\data"

# Remove the two obsolete/duplicate test cases MODEL_EXEC_13 and MODEL_EXEC_14
# (rows 19 and 20) - everything below shifts up automatically.
$ws.Rows("19:20").Delete() | Out-Null

# Move the active selection to where the removed rows used to be.
$ws.Range("A20").Select() | Out-Null
